$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-09 Tuesday", "2024-07-10 Wednesday"),
    @("281÷3=93, 2", "409÷9=45, 4"),
    @("295÷7=42, 1", "752÷6=125, 2"),
    @("568÷7=81, 1", "286÷7=40, 6"),
    @("866÷9=96, 2", "229÷8=28, 5"),
    @("377÷2=188, 1", "913÷7=130, 3"),
    @("477÷2=238, 1", "903÷2=451, 1"),
    @("922÷2=461, 0", "703÷4=175, 3"),
    @("824÷3=274, 2", "248÷5=49, 3"),
    @("709÷2=354, 1", "522÷4=130, 2"),
    @("144÷8=18, 0", "678÷7=96, 6"),
    @("642÷8=80, 2", "965÷8=120, 5"),
    @("668÷6=111, 2", "852÷5=170, 2"),
    @("369÷2=184, 1", "813÷2=406, 1"),
    @("836÷7=119, 3", "478÷6=79, 4"),
    @("930÷9=103, 3", "515÷4=128, 3"),
    @("949÷2=474, 1", "907÷2=453, 1"),
    @("889÷6=148, 1", "359÷9=39, 8"),
    @("430÷3=143, 1", "587÷9=65, 2"),
    @("616÷3=205, 1", "633÷8=79, 1"),
    @("198÷6=33, 0", "672÷6=112, 0"),
    @("752÷8=94, 0", "160÷9=17, 7"),
    @("305÷4=76, 1", "238÷7=34, 0"),
    @("719÷3=239, 2", "166÷4=41, 2"),
    @("514÷4=128, 2", "455÷2=227, 1"),
    @("713÷8=89, 1", "862÷3=287, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
